$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "37.889.87"
$ws.Range("E2").Value = "  +1.13%  "
$ws.Range("D3").Value = "2.038.40"
$ws.Range("E3").Value = "  +0.22%  "
$ws.Range("E4").Value = "  +0.04%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "227.43"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.99%  "
$ws.Range("E6").Value = "  -0.08%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "60.43"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +7.53%  "
$ws.Range("E8").Value = "  -0.05%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.383"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +0.36%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.0813"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +0.86%  "
$ws.Range("E11").Value = "  +0.74%  "
$ws.Range("E12").Value = "  +1.77%  "
$ws.Range("D13").Value = "2.338.75"
$ws.Range("E13").Value = "  -0.05%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "21.01"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +3.30%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.756"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +1.59%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "5.26"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +0.81%  "
$ws.Range("D17").Value = "2.064.51"
$ws.Range("E17").Value = "  +1.55%  "
$ws.Range("D18").Value = "37.848.62"
$ws.Range("E18").Value = "  +1.23%  "
$ws.Range("E19").Value = "  -1.82%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "69.70"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +0.93%  "
$ws.Range("D21").Value = "0.0₃0825"
$ws.Range("E21").Value = "  -0.19%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "224.00"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +0.20%  "
$ws.Range("E23").Value = "  +0.07%  "
$ws.Range("E24").Value = "  -0.92%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.19"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -2.70%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "166.47"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +1.09%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "9.16"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +0.20%  "
$ws.Range("E28").Value = "  -1.39%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "18.92"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +0.81%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.29"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -2.45%  "
$ws.Range("E31").Value = "  +2.02%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "4.46"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -0.57%  "
$ws.Range("E33").Value = "  +2.34%  "
$ws.Range("E34").Value = "  +0.55%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.0603"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -0.63%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "6.30"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +9.64%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.28"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -2.29%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "3.22"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -0.32%  "
$ws.Range("E39").Value = "  +0.03%  "
$ws.Range("D40").Value = "1.536.43"
$ws.Range("E40").Value = "  +4.12%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.0218"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +1.77%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "97.35"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +2.72%  "
$ws.Range("E43").Value = "  +0.61%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "16.59"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +1.81%  "
$ws.Range("E45").Value = "  -0.85%  "
$ws.Range("E46").Value = "  -0.17%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "4.00"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -5.39%  "
$ws.Range("E48").Value = "  +0.88%  "
$ws.Range("E49").Value = "  -0.47%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "7.04"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -0.98%  "
$ws.Range("D51").Value = "2.230.19"
$ws.Range("E51").Value = "  +0.17%  "
